$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 (header, "from" labels) - same shared-string slots 0-4 now hold new text
$ws.Range("B1").Value = "P_from_net1"
$ws.Range("C1").Value = "P_from_pv1"
$ws.Range("D1").Value = "P_from_bat1"
$ws.Range("E1").Value = "P_from_CHP1"
$ws.Range("F1").Value = "P_from_pvt1"

# Row 2 (values "to demand1")
$ws.Range("A2").Value = "P_to_demand1"
$ws.Range("B2").Value = "P_net1_demand1"
$ws.Range("C2").Value = "P_pv1_demand1"
$ws.Range("D2").Value = "P_bat1_demand1"
$ws.Range("E2").Value = "P_CHP1_demand1"
$ws.Range("F2").Value = "P_pvt1_demand1"

# Row 3 (values "to net1")
$ws.Range("A3").Value = "P_to_net1"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "P_pv1_net1"
$ws.Range("D3").Value = "P_bat1_net1"
$ws.Range("E3").Value = "P_CHP1_net1"
$ws.Range("F3").Value = "P_pvt1_net1"

# Row 4 (values "to bat1")
$ws.Range("A4").Value = "P_to_bat1"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "P_pv1_bat1"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "P_CHP1_bat1"
$ws.Range("F4").Value = "P_pvt1_bat1"
